$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in T4 (column E) grades for students (rows 2-7)
$ws.Range("E2").Value = 2.5
$ws.Range("E3").Value = 2.5
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0

# Update the active selection/view to reflect where the user ended up
$ws.Application.ActiveWindow.Zoom = 280
$ws.Range("E8").Select()
